# Update scripts with new TPM values (recomputed statistics).
# Target-cluster labels for rows 4/5/8/9 swap (MuSCs <-> Resolving-Mac)
# because of a reorder in the underlying cluster list; all other D-column
# values are unaffected. Numeric columns E-T are refreshed with the new
# TPM-derived statistics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("I2").Value = 0.9501423428580035
$ws.Range("J2").Value = 0.9501423428580037
$ws.Range("M2").Value = 509.3923236666667
$ws.Range("N2").Value = 1528.176971
$ws.Range("O2").Value = 0.831019558191033
$ws.Range("P2").Value = 0.8310195581910331
$ws.Range("Q2").Value = 21.42130558971311
$ws.Range("R2").Value = 192.791750307418
$ws.Range("S2").Value = 0.789586869980451
$ws.Range("T2").Value = 0.7895868699804514

# --- Row 3 ---
$ws.Range("I3").Value = 0.9501423428580035
$ws.Range("J3").Value = 0.9501423428580037
$ws.Range("N3").Value = 3.779073
$ws.Range("O3").Value = 0.002055052284145212
$ws.Range("P3").Value = 0.002055052284145212
$ws.Range("R3").Value = 0.476760291534
$ws.Range("S3").Value = 0.001952592191953423
$ws.Range("T3").Value = 0.001952592191953424

# --- Row 4 (Target cluster MuSCs -> still MuSCs, label list reordered) ---
$ws.Range("D4").Value = "MuSCs"
$ws.Range("I4").Value = 0.9501423428580035
$ws.Range("J4").Value = 0.9501423428580037
$ws.Range("M4").Value = 22.92703233333333
$ws.Range("N4").Value = 68.781097
$ws.Range("O4").Value = 0.03740302198339736
$ws.Range("P4").Value = 0.03740302198339737
$ws.Range("Q4").Value = 0.9641428483695554
$ws.Range("R4").Value = 8.677285635325999
$ws.Range("S4").Value = 0.03553819493727457
$ws.Range("T4").Value = 0.03553819493727459

# --- Row 5 (Target cluster Resolving-Mac -> still Resolving-Mac) ---
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("I5").Value = 0.9501423428580035
$ws.Range("J5").Value = 0.9501423428580037
$ws.Range("M5").Value = 79.39367866666667
$ws.Range("N5").Value = 238.181036
$ws.Range("O5").Value = 0.1295223675414243
$ws.Range("P5").Value = 0.1295223675414243
$ws.Range("Q5").Value = 3.338715904409778
$ws.Range("R5").Value = 30.048443139688
$ws.Range("S5").Value = 0.1230646857483243
$ws.Range("T5").Value = 0.1230646857483244

# --- Row 6 ---
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.002206666666666667
$ws.Range("H6").Value = 0.00662
$ws.Range("I6").Value = 0.04985765714199641
$ws.Range("J6").Value = 0.04985765714199642
$ws.Range("M6").Value = 509.3923236666667
$ws.Range("N6").Value = 1528.176971
$ws.Range("O6").Value = 0.831019558191033
$ws.Range("P6").Value = 0.8310195581910331
$ws.Range("Q6").Value = 1.124059060891111
$ws.Range("R6").Value = 10.11653154802
$ws.Range("S6").Value = 0.04143268821058186
$ws.Range("T6").Value = 0.04143268821058187

# --- Row 7 ---
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.002206666666666667
$ws.Range("H7").Value = 0.00662
$ws.Range("I7").Value = 0.04985765714199641
$ws.Range("J7").Value = 0.04985765714199642
$ws.Range("N7").Value = 3.779073
$ws.Range("O7").Value = 0.002055052284145212
$ws.Range("P7").Value = 0.002055052284145212
$ws.Range("Q7").Value = 0.00277971814
$ws.Range("R7").Value = 0.02501746326
$ws.Range("S7").Value = 0.0001024600921917886
$ws.Range("T7").Value = 0.0001024600921917886

# --- Row 8 (Target cluster MuSCs -> still MuSCs) ---
$ws.Range("D8").Value = "MuSCs"
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.002206666666666667
$ws.Range("H8").Value = 0.00662
$ws.Range("I8").Value = 0.04985765714199641
$ws.Range("J8").Value = 0.04985765714199642
$ws.Range("M8").Value = 22.92703233333333
$ws.Range("N8").Value = 68.781097
$ws.Range("O8").Value = 0.03740302198339736
$ws.Range("P8").Value = 0.03740302198339737
$ws.Range("Q8").Value = 0.05059231801555555
$ws.Range("R8").Value = 0.45533086214
$ws.Range("S8").Value = 0.00186482704612278
$ws.Range("T8").Value = 0.001864827046122781

# --- Row 9 (Target cluster Resolving-Mac -> still Resolving-Mac) ---
$ws.Range("D9").Value = "Resolving-Mac"
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.002206666666666667
$ws.Range("H9").Value = 0.00662
$ws.Range("I9").Value = 0.04985765714199641
$ws.Range("J9").Value = 0.04985765714199642
$ws.Range("M9").Value = 79.39367866666667
$ws.Range("N9").Value = 238.181036
$ws.Range("O9").Value = 0.1295223675414243
$ws.Range("P9").Value = 0.1295223675414243
$ws.Range("Q9").Value = 0.1751953842577778
$ws.Range("R9").Value = 1.57675845832
$ws.Range("S9").Value = 0.006457681793099979
$ws.Range("T9").Value = 0.00645768179309998
